$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (ECs) values
$ws.Range("G2").Value = 0.4859026666666666
$ws.Range("H2").Value = 1.457708
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.4581623333333333
$ws.Range("N2").Value = 1.374487
$ws.Range("O2").Value = 0.1060599910922654
$ws.Range("P2").Value = 0.1060599910922654
$ws.Range("Q2").Value = 0.2226222995328889
$ws.Range("R2").Value = 2.003600695796
$ws.Range("S2").Value = 0.1060599910922654
$ws.Range("T2").Value = 0.1060599910922654

# Update row 3 (FAPs) values
$ws.Range("G3").Value = 0.4859026666666666
$ws.Range("H3").Value = 1.457708
$ws.Range("M3").Value = 0.075101
$ws.Range("O3").Value = 0.01738512926863672
$ws.Range("P3").Value = 0.01738512926863672
$ws.Range("Q3").Value = 0.03649177616933334
$ws.Range("R3").Value = 0.328425985524
$ws.Range("S3").Value = 0.01738512926863672
$ws.Range("T3").Value = 0.01738512926863672

# Update row 4 - D4 changes to MuSCs; values change
$ws.Range("D4").Value = "MuSCs"
$ws.Range("G4").Value = 0.4859026666666666
$ws.Range("H4").Value = 1.457708
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.786578
$ws.Range("N4").Value = 11.359734
$ws.Range("O4").Value = 0.8765548796390978
$ws.Range("P4").Value = 0.8765548796390978
$ws.Range("Q4").Value = 1.839908347741333
$ws.Range("R4").Value = 16.559175129672
$ws.Range("S4").Value = 0.8765548796390978
$ws.Range("T4").Value = 0.8765548796390978

# Delete rows 5 and 6 (the extra MuSCs and Resolving-Mac rows)
$ws.Range("A5:T6").EntireRow.Delete()
